$wb = $excel.ActiveWorkbook

# New handoff/handback report identifiers for this run
# (replaces 882fc755-095b-4cc6-af4b-658cf8c09ce0 and 9b8ec366-ef01-4bbb-b031-849c0146b210)
$newUuid1 = "f5586bb4-522f-481d-91e0-a6c6ccdb70b5"
$newUuid2 = "ffff521b305b-8a56-4183-b1b0-1aaedbf75463"
$newHash  = "7e9b4f69c0ac9bd11e653f5ef308e0be86c4bcf0"

$newMd1 = "$newUuid1.md"
$newMd2 = "$newUuid2.md"
$newZhCnXlf = "$newUuid1.$newHash.zh-cn.xlf"
$newDeDeXlf = "$newUuid1.$newHash.de-de.xlf"

$zhCnHandoffDt  = "2016-03-24 01:10:43"
$zhCnHandbackDt = "2016-03-24 01:11:06"
$deDeHandoffDt  = "2016-03-24 01:10:48"
$deDeHandbackDt = "2016-03-24 01:11:15"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = $newMd1 }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = $newMd2 }
}
$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = $newMd1 }
    if ($addr -eq '$D$2') { $hl.TextToDisplay = $newZhCnXlf }
    if ($addr -eq '$F$2') { $hl.TextToDisplay = $newMd1 }
    if ($addr -eq '$G$2') { $hl.TextToDisplay = $newZhCnXlf }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = $newMd2 }
    if ($addr -eq '$D$3') { $hl.TextToDisplay = $newZhCnXlf }
    if ($addr -eq '$F$3') { $hl.TextToDisplay = $newMd2 }
    if ($addr -eq '$G$3') { $hl.TextToDisplay = $newZhCnXlf }
}
$wsZhCn.Range("A2").Value = $newMd1
$wsZhCn.Range("D2").Value = $newZhCnXlf
$wsZhCn.Range("E2").Value = $zhCnHandoffDt
$wsZhCn.Range("F2").Value = $newMd1
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $zhCnHandbackDt
$wsZhCn.Range("A3").Value = $newMd2
$wsZhCn.Range("D3").Value = $newZhCnXlf
$wsZhCn.Range("E3").Value = $zhCnHandoffDt
$wsZhCn.Range("F3").Value = $newMd2
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = $zhCnHandbackDt

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = $newMd1 }
    if ($addr -eq '$D$2') { $hl.TextToDisplay = $newDeDeXlf }
    if ($addr -eq '$F$2') { $hl.TextToDisplay = $newMd1 }
    if ($addr -eq '$G$2') { $hl.TextToDisplay = $newDeDeXlf }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = $newMd2 }
    if ($addr -eq '$D$3') { $hl.TextToDisplay = $newDeDeXlf }
    if ($addr -eq '$F$3') { $hl.TextToDisplay = $newMd2 }
    if ($addr -eq '$G$3') { $hl.TextToDisplay = $newDeDeXlf }
}
$wsDeDe.Range("A2").Value = $newMd1
$wsDeDe.Range("D2").Value = $newDeDeXlf
$wsDeDe.Range("E2").Value = $deDeHandoffDt
$wsDeDe.Range("F2").Value = $newMd1
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $deDeHandbackDt
$wsDeDe.Range("A3").Value = $newMd2
$wsDeDe.Range("D3").Value = $newDeDeXlf
$wsDeDe.Range("E3").Value = $deDeHandoffDt
$wsDeDe.Range("F3").Value = $newMd2
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = $deDeHandbackDt
